$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.897.20'
$ws.Range("E2").Value = '  -1.06%  '
$ws.Range("D3").Value = '1.641.21'
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("E4").Value = '  -0.58%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.83'
$ws.Range("E5").Value = '  -0.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5043'
$ws.Range("E6").Value = '  -1.70%  '
$ws.Range("E7").Value = '  -0.70%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2579'
$ws.Range("E8").Value = '  -0.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06386'
$ws.Range("E9").Value = '  -0.92%  '
$ws.Range("E10").Value = '  -1.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07743'
$ws.Range("E11").Value = '  -0.91%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.267'
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '1.866.95'
$ws.Range("E13").Value = '  -0.95%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.600.01'
$ws.Range("E14").Value = '  -3.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5456'
$ws.Range("E15").Value = '  -1.06%  '
$ws.Range("D16").Value = '0.0₅7901'
$ws.Range("E16").Value = '  -1.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.20'
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").Value = '25.916.48'
$ws.Range("E18").Value = '  -1.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.003'
$ws.Range("E19").Value = '  -0.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '202.36'
$ws.Range("E20").Value = '  -2.99%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.392'
$ws.Range("E21").Value = '  -0.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.890'
$ws.Range("E22").Value = '  -1.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.978'
$ws.Range("E23").Value = '  -1.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.004'
$ws.Range("E24").Value = '  -0.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.861'
$ws.Range("E25").Value = '  +3.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '140.78'
$ws.Range("E26").Value = '  -2.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1136'
$ws.Range("E27").Value = '  -3.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.65'
$ws.Range("E28").Value = '  -1.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.780'
$ws.Range("E29").Value = '  -2.78%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.245'
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04982'
$ws.Range("E31").Value = '  -1.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.274'
$ws.Range("E32").Value = '  -2.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.201'
$ws.Range("E33").Value = '  -1.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.548'
$ws.Range("E34").Value = '  -0.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.368'
$ws.Range("E35").Value = '  +0.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.630'
$ws.Range("E36").Value = '  -3.97%  '
$ws.Range("E37").Value = '  -3.18%  '
$ws.Range("D38").Value = '1.151.63'
$ws.Range("E38").Value = '  -1.96%  '
$ws.Range("E39").Value = '  -1.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01564'
$ws.Range("E40").Value = '  -1.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.003'
$ws.Range("E41").Value = '  -0.69%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.687'
$ws.Range("E42").Value = '  +0.46%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8079'
$ws.Range("E43").Value = '  -2.15%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.80'
$ws.Range("E44").Value = '  -0.62%  '
$ws.Range("D45").Value = '1.777.89'
$ws.Range("E45").Value = '  -1.11%  '
$ws.Range("E46").Value = '  +2.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4537'
$ws.Range("E47").Value = '  -0.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.008'
$ws.Range("E48").Value = '  +0.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.82'
$ws.Range("E49").Value = '  -1.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05055'
$ws.Range("E50").Value = '  -0.62%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.001'
$ws.Range("E51").Value = '  -1.07%  '
